$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Contest 38 (row 47, "KXI vs DC") results entered ---
$ws.Range("E47").Value = 80
$ws.Range("H47").Value = 20
$ws.Range("K47").Value = 40
$ws.Range("N47").Value = 100
$ws.Range("Q47").Value = 60
$ws.Range("T47").Value = 0

# --- Insert a new blank "template" row at 56, pushing the old row 56
#     (and everything below it) down by one. This keeps the SUM(...)
#     ranges on the totals row correctly auto-expanded to include it. ---
$ws.Rows("56:56").Insert(-4121)

# The freshly inserted row 56 does not inherit formatting/formulas the
# way a real Excel "insert" with CopyOrigin would, so restore that by
# pasting the (now shifted) row 57's formatting back onto row 56 --
# row 57 is exactly what row 56 used to look like before the insert.
$ws.Range("A57:E57").Copy()
$ws.Range("A56:E56").PasteSpecial(-4122)
$ws.Range("G57:H57").Copy()
$ws.Range("G56:H56").PasteSpecial(-4122)
$ws.Range("J57:K57").Copy()
$ws.Range("J56:K56").PasteSpecial(-4122)
$ws.Range("M57:N57").Copy()
$ws.Range("M56:N56").PasteSpecial(-4122)
$ws.Range("P57:Q57").Copy()
$ws.Range("P56:Q56").PasteSpecial(-4122)
$ws.Range("S57:T57").Copy()
$ws.Range("S56:T56").PasteSpecial(-4122)

# --- Contest 46 (new row 56, "SRH vs DC") ---
$ws.Range("A56").Value = 46
$ws.Range("B56").Value = 2
$ws.Range("C56").Value = "SRH vs DC"

$ws.Range("D56").Formula = '=IF(ISERROR(VLOOKUP(RANK(E56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(E56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("G56").Formula = '=IF(ISERROR(VLOOKUP(RANK(H56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(H56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("J56").Formula = '=IF(ISERROR(VLOOKUP(RANK(K56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(K56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("M56").Formula = '=IF(ISERROR(VLOOKUP(RANK(N56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(N56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("P56").Formula = '=IF(ISERROR(VLOOKUP(RANK(Q56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(Q56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("S56").Formula = '=IF(ISERROR(VLOOKUP(RANK(T56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(T56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'

# --- Keep the selection in sync with where Excel would land (one row
#     further down than before, since a row was inserted above it) ---
$ws.Range("U61").Select()
